$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C,D,E,F,G,I,J,K,L,N (rows 2-25); B,H,M,O unchanged
$cols = @(3,4,5,6,7,9,10,11,12,14)

$data = @(
    @(0.06488711199025943, 0.01635300175438559, 0.1098494444106812, 4.571276472368226, 0.002653899182221333, 3.332893248685494, 0.2579758678232764, 4.640079425027864, 0.2137478850433965, 2.294238401452205),
    @(0.06387223968761191, 0.01589063522901313, 0.1094484679372911, 4.570580968545528, 0.002659709070497144, 3.332587146422043, 0.2580690294620993, 4.495262951214954, 0.2136980212008766, 2.321807395555494),
    @(0.06326923176878552, 0.01560775708728457, 0.1092418690032009, 4.572543650682661, 0.002663463562251276, 3.334155201672814, 0.2582351773818132, 4.408773118721513, 0.2137476198980011, 2.339536884406909),
    @(0.0630285890193818, 0.01549274928614253, 0.1091676517623199, 4.573943089524363, 0.002665040783951232, 3.335234776617057, 0.2583302681707309, 4.374137341501751, 0.2137880153359859, 2.346963306593604),
    @(0.06298893856074272, 0.01547366882466861, 0.1091559308776962, 4.574211644326311, 0.002665305538216423, 3.335440620786059, 0.2583477115542721, 4.368422882636423, 0.2137959423982068, 2.348208624536733),
    @(0.06326596573953225, 0.01560620495330411, 0.1092408276803329, 4.572560097858215, 0.002663484641470579, 3.334167978629281, 0.2582363489460491, 4.408303541895236, 0.2137480829453438, 2.339636224055258),
    @(0.06453302188368326, 0.01619337349302086, 0.109702976923014, 4.570539883062537, 0.00265586367786506, 3.332422641837397, 0.2579853705039241, 4.589642069682839, 0.2137140502910171, 2.303577636791626),
    @(0.06717644597885197, 0.01735243875618409, 0.1109229765986726, 4.585604076744289, 0.002642396825735638, 3.34298210647502, 0.2583584615253898, 4.964584753166832, 0.2142834672314393, 2.239236668398085),
    @(0.06921431353443097, 0.01820817593381818, 0.1120101446443655, 4.608368987751589, 0.002633393158765207, 3.359338603307251, 0.2591616429042745, 5.25198340220021, 0.2150895720293988, 2.195855542515956),
    @(0.07016199142648816, 0.01859828948085251, 0.1125460990398075, 4.621287627682534, 0.002629488262557946, 3.368663556006624, 0.2596423071166925, 5.385351927123793, 0.2155405142176647, 2.176966661962576),
    @(0.07052379995145941, 0.01874612639643658, 0.1127549948032538, 4.626549828503727, 0.002628036861999286, 3.37246694020277, 0.2598409285372512, 5.436235286365502, 0.2157233859961281, 2.169935713879372),
    @(0.07044574748825028, 0.01871428237212314, 0.1127097412881817, 4.625400028037689, 0.002628348235613448, 3.371635685538138, 0.2597974129612126, 5.425259737625026, 0.2156834627550097, 2.171444529781118),
    @(0.07019169877911224, 0.01861044998074846, 0.112563166006673, 4.621713123223486, 0.002629368308757461, 3.368970999532294, 0.2596583148953187, 5.389530515346337, 0.2155553165103257, 2.176385777647978),
    @(0.07003646921126006, 0.01854686357905422, 0.1124741578440229, 4.619503047659578, 0.002629996683632723, 3.367374292438498, 0.2595752764609642, 5.367694814539732, 0.2154784001814321, 2.179428313383628),
    @(0.0691527937884473, 0.01818269697274033, 0.111975950786956, 4.607576435598133, 0.002633652181658759, 3.358767220367014, 0.2591325524899446, 5.243320491632574, 0.2150617966787749, 2.197107011836176),
    @(0.06861595657770891, 0.01795949860589019, 0.1116809112707067, 4.60091740848452, 0.002635943498982074, 3.353970568203692, 0.2588905025448511, 5.167695368072259, 0.2148277969259524, 2.208169055235015),
    @(0.06830912635540898, 0.01783119995834426, 0.1115151093903748, 4.597328430180056, 0.002637279383542074, 3.351388942735824, 0.258762132140447, 5.124445199322622, 0.2147011358069264, 2.214611206051268),
    @(0.06820557362210167, 0.0177877741652992, 0.1114596413959532, 4.596154625699199, 0.002637734783550021, 3.350545258853003, 0.2587205309310789, 5.109843882817017, 0.214659612539073, 2.216806062707377),
    @(0.06867290276400695, 0.01798325033934134, 0.1117119154847686, 4.601601306520422, 0.002635697724490005, 3.35446282291278, 0.2589151460288548, 5.175720175884862, 0.2148518859738928, 2.206983246004999),
    @(0.07026623935394127, 0.01864094518543169, 0.1126060575573256, 4.622785996148139, 0.002629067948903913, 3.36974628401606, 0.2596987205237653, 5.400014739786627, 0.2155926275405946, 2.174931103719224),
    @(0.07132472262698286, 0.01907141968062831, 0.1132250576676803, 4.638789978877469, 0.002624894054505401, 3.381322263010162, 0.2603076271762816, 5.548817503358691, 0.2161473287566054, 2.154693618059195),
    @(0.07075822939847853, 0.01884161283970442, 0.1128915209011261, 4.630050282667042, 0.002627107238812811, 3.374998276806977, 0.2599737767100621, 5.469195642051204, 0.2158448169254825, 2.165429640352064),
    @(0.06864715177048453, 0.01797251211524653, 0.1116978865795133, 4.601291370713085, 0.002635808781418607, 3.354239726348425, 0.2589039711116783, 5.172091451079154, 0.2148409708065699, 2.207519093614227),
    @(0.0664444441799219, 0.01703811498058982, 0.1105593933264473, 4.579481894849266, 0.002645882841422019, 3.338620990348659, 0.2581647394871496, 4.861068388581259, 0.2137140502910171, 2.255959989681939)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($r + 2, $cols[$c]).Value2 = $data[$r][$c]
    }
}
